# Update "想去人数" (want-to-go count) for a handful of events.
# Same underlying data is duplicated on the "展览" sheet and the
# "全部类型" sheet, so both need the refreshed counts.
$wb = $excel.ActiveWorkbook

$updates = @{
    "F5"  = 4634
    "F9"  = 913
    "F11" = 1074
    "F15" = 13
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
